$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos list refresh: updated Price (D) and Volume(1h) (E) figures for each
# coin row, plus three rows whose rank order swapped (Hedera/EnergySwap and
# Filecoin/FirstDigitalUSD), which moves their Coin name (B) and Link (C) too.
# Price values are written with a leading apostrophe so Excel stores them as
# literal text (matching the source data) instead of auto-converting look-alike
# numbers such as "1.00" or "0.437" into numeric values.

# Row 2
$ws.Range("D2").Value = '''57.376.50'
$ws.Range("E2").Value = '  +1.40%  '

# Row 3
$ws.Range("D3").Value = '''3.004.64'
$ws.Range("E3").Value = '  -0.01%  '

# Row 4
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  +0.02%  '

# Row 5
$ws.Range("D5").Value = '''509.21'
$ws.Range("E5").Value = '  +0.26%  '

# Row 6
$ws.Range("D6").Value = '''139.08'
$ws.Range("E6").Value = '  +1.41%  '

# Row 7
$ws.Range("D7").Value = '''1.00'
$ws.Range("E7").Value = '  -0.02%  '

# Row 8
$ws.Range("D8").Value = '''0.437'
$ws.Range("E8").Value = '  +0.84%  '

# Row 9
$ws.Range("D9").Value = '''7.52'
$ws.Range("E9").Value = '  -1.03%  '

# Row 10
$ws.Range("E10").Value = '  +1.63%  '

# Row 11
$ws.Range("D11").Value = '''0.366'
$ws.Range("E11").Value = '  +3.39%  '

# Row 12
$ws.Range("D12").Value = '''3.521.81'
$ws.Range("E12").Value = '  +0.03%  '

# Row 13
$ws.Range("E13").Value = '  +1.32%  '

# Row 14
$ws.Range("D14").Value = '''26.44'
$ws.Range("E14").Value = '  +3.30%  '

# Row 15
$ws.Range("D15").Value = '''0.0000164'
$ws.Range("E15").Value = '  +6.50%  '

# Row 16
$ws.Range("D16").Value = '''57.380.65'
$ws.Range("E16").Value = '  +1.38%  '

# Row 17
$ws.Range("D17").Value = '''6.22'
$ws.Range("E17").Value = '  +6.72%  '

# Row 18
$ws.Range("D18").Value = '''3.010.18'
$ws.Range("E18").Value = '  +0.16%  '

# Row 19
$ws.Range("D19").Value = '''12.78'
$ws.Range("E19").Value = '  +2.26%  '

# Row 20
$ws.Range("D20").Value = '''7.95'
$ws.Range("E20").Value = '  +1.22%  '

# Row 21
$ws.Range("D21").Value = '''328.99'
$ws.Range("E21").Value = '  +0.61%  '

# Row 22
$ws.Range("D22").Value = '''1.00'
$ws.Range("E22").Value = '  +0.04%  '

# Row 23
$ws.Range("D23").Value = '''0.497'
$ws.Range("E23").Value = '  +3.82%  '

# Row 24
$ws.Range("D24").Value = '''64.26'
$ws.Range("E24").Value = '  +2.80%  '

# Row 25
$ws.Range("E25").Value = '  +0.98%  '

# Row 26
$ws.Range("E26").Value = '  -0.26%  '

# Row 27
$ws.Range("D27").Value = '''0.0₃0915'
$ws.Range("E27").Value = '  -0.61%  '

# Row 28
$ws.Range("D28").Value = '''6.78'
$ws.Range("E28").Value = '  +3.14%  '

# Row 29
$ws.Range("D29").Value = '''7.45'
$ws.Range("E29").Value = '  +5.95%  '

# Row 30
$ws.Range("D30").Value = '''1.81'
$ws.Range("E30").Value = '  +2.31%  '

# Row 31
$ws.Range("D31").Value = '''1.19'
$ws.Range("E31").Value = '  -4.46%  '

# Row 32
$ws.Range("D32").Value = '''20.62'
$ws.Range("E32").Value = '  -0.25%  '

# Row 33
$ws.Range("D33").Value = '''4.71'
$ws.Range("E33").Value = '  +4.16%  '

# Row 34
$ws.Range("D34").Value = '''154.23'
$ws.Range("E34").Value = '  -1.53%  '

# Row 35
$ws.Range("D35").Value = '''5.87'
$ws.Range("E35").Value = '  +3.65%  '

# Row 36
$ws.Range("D36").Value = '''1.28'
$ws.Range("E36").Value = '  +0.41%  '

# Row 37
$ws.Range("B37").Value = 'EnergySwap'
$ws.Range("C37").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D37").Value = '''24.47'
$ws.Range("E37").Value = '  +2.56%  '

# Row 38
$ws.Range("B38").Value = 'Hedera'
$ws.Range("C38").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D38").Value = '''0.0680'
$ws.Range("E38").Value = '  +0.75%  '

# Row 39
$ws.Range("D39").Value = '''3.039.60'
$ws.Range("E39").Value = '  +0.03%  '

# Row 40
$ws.Range("E40").Value = '  +1.57%  '

# Row 41
$ws.Range("B41").Value = 'Filecoin'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D41").Value = '''3.84'
$ws.Range("E41").Value = '  +6.21%  '

# Row 42
$ws.Range("B42").Value = 'FirstDigitalUSD'
$ws.Range("C42").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D42").Value = '''1.00'
$ws.Range("E42").Value = '  +0.02%  '

# Row 43
$ws.Range("D43").Value = '''2.295.21'
$ws.Range("E43").Value = '  +1.23%  '

# Row 44
$ws.Range("E44").Value = '  -0.02%  '

# Row 45
$ws.Range("D45").Value = '''1.41'
$ws.Range("E45").Value = '  +0.15%  '

# Row 46
$ws.Range("D46").Value = '''0.986'
$ws.Range("E46").Value = '  -1.46%  '

# Row 47
$ws.Range("D47").Value = '''6.02'
$ws.Range("E47").Value = '  +3.91%  '

# Row 48
$ws.Range("D48").Value = '''0.0239'
$ws.Range("E48").Value = '  +0.93%  '

# Row 49
$ws.Range("D49").Value = '''19.44'
$ws.Range("E49").Value = '  +1.23%  '

# Row 50
$ws.Range("D50").Value = '''1.84'
$ws.Range("E50").Value = '  -7.09%  '

# Row 51
$ws.Range("D51").Value = '''0.0893'
$ws.Range("E51").Value = '  +1.96%  '
